$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in E1
$ws.Range("E1").Value = "MP"

# Replace each E column value (rows 2-183) with its reciprocal (1/x)
for ($r = 2; $r -le 183; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value2
    if ($old -ne $null -and $old -ne 0) {
        $cell.Value2 = 1 / $old
    }
}
